{"js": "// Add \", Oracle 12c\" before \", Neo4j and MongoDB\" in the \"Technology/Tools Uses\" line,\n// add \", Big Data Integration and Processing\" before \", Big Data Essentials HDFS MapReduce\n// and Spark RDD.\" in the Specialization list, and relocate the (single) \"_GoBack\" bookmark\n// from its old spot (inside \"Gold Badge for Python\") to right after the newly inserted\n// \"Big Data Integration and Processing\" text.\n\nconst body = context.document.body;\n\n// 1) Remove the existing \"_GoBack\" bookmark first so re-inserting it below actually\n//    moves it (Word bookmark names are unique; deleting first then inserting relocates it).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) \", Oracle 12c\" goes right before \", Neo4j and MongoDB\".\nconst neoResults = body.search(\", Neo4j and MongoDB\", { matchCase: true, matchWholeWord: false });\nneoResults.load(\"items\");\nawait context.sync();\n\nif (neoResults.items.length > 0) {\n  const neoStart = neoResults.items[0].getRange(\"Start\");\n  neoStart.insertText(\", Oracle 12c\", \"Before\");\n  await context.sync();\n}\n\n// 3) \", Big Data Integration and Processing\" goes right before\n//    \", Big Data Essentials HDFS MapReduce and Spark RDD.\"\nconst bdeResults = body.search(\", Big Data Essentials HDFS MapReduce and Spark RDD.\", { matchCase: true, matchWholeWord: false });\nbdeResults.load(\"items\");\nawait context.sync();\n\nif (bdeResults.items.length > 0) {\n  const bdeStart = bdeResults.items[0].getRange(\"Start\");\n  bdeStart.insertText(\", Big Data Integration and Processing\", \"Before\");\n  await context.sync();\n}\n\n// 4) Re-insert the \"_GoBack\" bookmark immediately before \", Big Data Essentials...\" \u2014\n//    i.e. right after the text we just inserted.\nconst bdeResults2 = body.search(\", Big Data Essentials HDFS MapReduce and Spark RDD.\", { matchCase: true, matchWholeWord: false });\nbdeResults2.load(\"items\");\nawait context.sync();\n\nif (bdeResults2.items.length > 0) {\n  const bdeStart2 = bdeResults2.items[0].getRange(\"Start\");\n  bdeStart2.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Add \", Oracle 12c\" before \", Neo4j and MongoDB\" in the \"Technology/Tools Uses\" line,\n# add \", Big Data Integration and Processing\" before \", Big Data Essentials HDFS MapReduce\n# and Spark RDD.\" in the Specialization list, and relocate the (single) \"_GoBack\" bookmark\n# from its old spot (inside \"Gold Badge for Python\") to right after the newly inserted\n# \"Big Data Integration and Processing\" text.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the existing \"_GoBack\" bookmark first so re-adding it below actually moves it\n#    (bookmark names are unique; delete-then-add relocates it to the new range).\n$oldBookmark = $d.Bookmarks.Item(\"_GoBack\")\n$oldBookmark.Delete()\n\n# 2) \", Oracle 12c\" goes right before \", Neo4j and MongoDB\".\n$r1 = $d.Content\n$r1.Find.Execute(\", Neo4j and MongoDB\") | Out-Null\n$r1.Collapse(1)            # wdCollapseStart\n$r1.InsertBefore(\", Oracle 12c\")\n\n# 3) \", Big Data Integration and Processing\" goes right before\n#    \", Big Data Essentials HDFS MapReduce and Spark RDD.\"\n$r2 = $d.Content\n$r2.Find.Execute(\", Big Data Essentials HDFS MapReduce and Spark RDD.\") | Out-Null\n$r2.Collapse(1)            # wdCollapseStart\n$r2.InsertBefore(\", Big Data Integration and Processing\")\n\n# 4) Re-add the \"_GoBack\" bookmark immediately before \", Big Data Essentials...\" \u2014\n#    i.e. right after the text we just inserted.\n$r3 = $d.Content\n$r3.Find.Execute(\", Big Data Essentials HDFS MapReduce and Spark RDD.\") | Out-Null\n$r3.Collapse(1)            # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $r3)\n"}
